# Add "Sheet2" as a new worksheet, placed right after Sheet1 (i.e. at the end)
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Populate Sheet2 content
$ws2.Range("A1").Value = "Lampiran Pengajuan SK"

$ws2.Range("A3").Value = "no"
$ws2.Range("B3").Value = "Jabatan"
$ws2.Range("C3").Value = "Nama"
$ws2.Range("D3").Value = "Uraian Tugas"

$ws2.Range("A13").Value = "Lampiran Honorarium (Jika dihonorkan)"

$ws2.Range("A15").Value = "Kode Akun"
$ws2.Range("B15").Value = "Jabatan"
$ws2.Range("C15").Value = "Nominal"

# Select A15 as the active cell on Sheet2
$ws2.Range("A15").Select()

# Make Sheet2 the active sheet
$ws2.Activate()

# Sheet1's selection changes to C4 (no longer the tab-selected sheet)
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("C4").Select()

# Re-activate Sheet2 to ensure it's the final active tab
$ws2.Activate()
